# Updates cryptos list values (Price + Volume(1h) columns) to reflect the
# latest scrape, per the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cells needing forced text number format (numeric-looking values)
$textFormatCells = @("D5", "D6", "D9", "D16", "D18", "D21", "D24", "D29", "D31", "D36", "D37", "D43", "D47", "D51")
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# apply new values
$ws.Range("D2").Value = '64.433.58'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '3.517.39'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '592.76'
$ws.Range("E5").Value = '  +1.55%  '
$ws.Range("D6").Value = '134.79'
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D9").Value = '7.59'
$ws.Range("E9").Value = '  +6.44%  '
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("E11").Value = '  +4.20%  '
$ws.Range("D12").Value = '4.115.53'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("E14").Value = '  +1.03%  '
$ws.Range("D15").Value = '3.518.43'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '25.97'
$ws.Range("D17").Value = '64.418.76'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '9.95'
$ws.Range("E18").Value = '  +2.11%  '
$ws.Range("E19").Value = '  +3.52%  '
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D21").Value = '394.82'
$ws.Range("E21").Value = '  +2.83%  '
$ws.Range("E22").Value = '  +1.47%  '
$ws.Range("D23").Value = '3.657.52'
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").Value = '74.80'
$ws.Range("E24").Value = '  +1.29%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("E27").Value = '  +2.67%  '
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").Value = '  +2.50%  '
$ws.Range("D31").Value = '8.34'
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("E32").Value = '  -5.60%  '
$ws.Range("E33").Value = '  +7.55%  '
$ws.Range("D34").Value = '3.546.29'
$ws.Range("E34").Value = '  +0.66%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = '23.48'
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("D37").Value = '5.37'
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("E39").Value = '  +1.74%  '
$ws.Range("E40").Value = '  +1.75%  '
$ws.Range("E41").Value = '  +1.64%  '
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").Value = '25.48'
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("E46").Value = '  +3.46%  '
$ws.Range("D47").Value = '1.18'
$ws.Range("E47").Value = '  -2.82%  '
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("D49").Value = '2.412.20'
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("D51").Value = '0.0260'
$ws.Range("E51").Value = '  +0.22%  '

Write-Output "Updated cryptos list values."
